$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.496.60'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '2.484.83'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '527.03'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').Value = '133.86'
$ws.Range('E6').Value = '  -3.07%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = '0.100'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = '0.156'
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('D11').Value = '5.37'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').Value = '0.342'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').Value = '2.925.78'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = '58.426.69'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').Value = '22.43'
$ws.Range('E15').Value = '  -3.25%  '
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Value = '2.477.46'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '10.91'
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('D20').Value = '320.55'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = '5.81'
$ws.Range('E22').Value = '  -1.33%  '
$ws.Range('D23').Value = '64.34'
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('D24').Value = '0.415'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = '0.162'
$ws.Range('E25').Value = '  -2.09%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = '7.45'
$ws.Range('E27').Value = '  -2.70%  '
$ws.Range('D28').Value = '0.0₃0752'
$ws.Range('E28').Value = '  -3.49%  '
$ws.Range('D29').Value = '6.45'
$ws.Range('E29').Value = '  -3.96%  '
$ws.Range('D30').Value = '1.72'
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('D31').Value = '167.14'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').Value = '1.14'
$ws.Range('E32').Value = '  -5.45%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').Value = '18.25'
$ws.Range('E35').Value = '  -1.60%  '
$ws.Range('D36').Value = '1.34'
$ws.Range('E36').Value = '  -9.56%  '
$ws.Range('D37').Value = '3.98'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('E38').Value = '  -4.41%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').Value = '0.793'
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.53'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('D41').Value = '276.53'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('D42').Value = '4.94'
$ws.Range('E42').Value = '  -5.86%  '
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D44').Value = '127.93'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('D45').Value = '0.0914'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('D47').Value = '0.0216'
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('D48').Value = '17.17'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('D49').Value = '1.736.09'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('D50').Value = '0.979'
$ws.Range('E50').Value = '  -1.23%  '
$ws.Range('D51').Value = '4.71'
$ws.Range('E51').Value = '  -1.49%  '
